# Auto-applied update: dades i banners 2026-02-16 21:20
# Updates DATA_EXTRACCIO timestamps and several measurement cells
# in the meteocat daily summary sheet to match the latest extraction run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues
$xlPasteValues = -4163

$ws.Range("E2").Value = "2026-02-16 21:18:34"
$ws.Range("E3").Value = "2026-02-16 21:18:37"
$ws.Range("N3").Value = "-2.5 °C 20:59 TU"
$ws.Range("O3").Value = "-1.0 °C"
$ws.Range("E4").Value = "2026-02-16 21:18:39"
$ws.Range("O4").Value = "13.5 °C"
$ws.Range("E5").Value = "2026-02-16 21:18:41"
$ws.Range("I5").Value = "25.1 mm"
$ws.Range("L5").Value = "54.7 km/h - 326º 20:39 TU"
$ws.Range("N5").Value = "-3.3 °C 20:59 TU"
$ws.Range("E6").Value = "2026-02-16 21:18:44"
$ws.Range("J6").Value = "1012.4 hPa"
$ws.Range("E7").Value = "2026-02-16 21:18:46"
$ws.Range("E8").Value = "2026-02-16 21:18:49"
$ws.Range("E9").Value = "2026-02-16 21:18:51"
$ws.Range("L9").Value = "49.3 km/h - 335º 20:54 TU"
$ws.Range("O9").Value = "11.5 °C"
$ws.Range("E10").Value = "2026-02-16 21:18:54"
$ws.Range("O10").Value = "10.8 °C"
$ws.Range("E11").Value = "2026-02-16 21:18:56"
$ws.Range("O11").Value = "7.0 °C"
$ws.Range("E12").Value = "2026-02-16 21:18:59"
$ws.Range("E13").Value = "2026-02-16 21:19:01"
$ws.Range("E14").Value = "2026-02-16 21:19:04"
$ws.Range("K14").Value = "9.8 MJ/m2"
$ws.Range("E15").Value = "2026-02-16 21:19:06"
$ws.Range("O15").Value = "11.8 °C"
$ws.Range("E16").Value = "2026-02-16 21:19:09"
$ws.Range("L16").Value = "131.8 km/h - 171º 20:51 TU"
$ws.Range("N16").Value = "-2.5 °C 20:59 TU"
$ws.Range("E17").Value = "2026-02-16 21:19:11"
$ws.Range("E18").Value = "2026-02-16 21:19:14"
$ws.Range("E19").Value = "2026-02-16 21:19:16"
$ws.Range("E20").Value = "2026-02-16 21:19:19"
$ws.Range("O20").Value = "-0.7 °C"
$ws.Range("E21").Value = "2026-02-16 21:19:21"
$ws.Range("O21").Value = "8.8 °C"
$ws.Range("E22").Value = "2026-02-16 21:19:24"
$ws.Range("E23").Value = "2026-02-16 21:19:26"
$ws.Range("I23").Value = "15.8 mm"
$ws.Range("N23").Value = "-3.7 °C 20:59 TU"
$ws.Range("E24").Value = "2026-02-16 21:19:29"
$ws.Range("E25").Value = "2026-02-16 21:19:31"
$ws.Range("I25").Value = "6.5 mm"
$ws.Range("N25").Value = "-1.1 °C 20:58 TU"
$ws.Range("O25").Value = "0.6 °C"
$ws.Range("E26").Value = "2026-02-16 21:19:34"
$ws.Range("E27").Value = "2026-02-16 21:19:36"
$ws.Range("E28").Value = "2026-02-16 21:19:38"
$ws.Range("O28").Value = "9.6 °C"
$ws.Range("E29").Value = "2026-02-16 21:19:41"
$ws.Range("E30").Value = "2026-02-16 21:19:43"
$ws.Range("L30").Value = "36.7 km/h - 33º 20:58 TU"
$ws.Range("O30").Value = "11.8 °C"
$ws.Range("E31").Value = "2026-02-16 21:19:46"
$ws.Range("N31").Value = "11.5 °C 20:58 TU"
$ws.Range("E32").Value = "2026-02-16 21:19:48"
$ws.Range("L32").Value = "54.0 km/h - 282º 20:56 TU"
$ws.Range("E33").Value = "2026-02-16 21:19:51"
$ws.Range("E34").Value = "2026-02-16 21:19:54"
$ws.Range("I34").Value = "3.6 mm"
$ws.Range("N34").Value = "1.8 °C 20:51 TU"
$ws.Range("E35").Value = "2026-02-16 21:19:56"
$ws.Range("I35").Value = "1.8 mm"
$ws.Range("E36").Value = "2026-02-16 21:19:58"
$ws.Range("L36").Value = "65.9 km/h - 339º 20:42 TU"
$ws.Range("O36").Value = "12.1 °C"
$ws.Range("E37").Value = "2026-02-16 21:20:01"
$ws.Range("E38").Value = "2026-02-16 21:20:03"
$ws.Range("O38").Value = "11.9 °C"
$ws.Range("E39").Value = "2026-02-16 21:20:06"
$ws.Range("I39").Value = "4.4 mm"
$ws.Range("N39").Value = "-2.4 °C 20:59 TU"
$ws.Range("O39").Value = "0.2 °C"
$ws.Range("E40").Value = "2026-02-16 21:20:08"
$ws.Range("J40").Value = "1016.2 hPa"
$ws.Range("E41").Value = "2026-02-16 21:20:11"
$ws.Range("J41").Value = "1014.7 hPa"
$ws.Range("E42").Value = "2026-02-16 21:20:14"
$ws.Range("E43").Value = "2026-02-16 21:20:16"
$ws.Range("O43").Value = "8.9 °C"
$ws.Range("E44").Value = "2026-02-16 21:20:18"
$ws.Range("I44").Value = "12.7 mm"
$ws.Range("L44").Value = "69.8 km/h - 63º 20:44 TU"
$ws.Range("N44").Value = "-2.4 °C 20:59 TU"
$ws.Range("E45").Value = "2026-02-16 21:20:21"
$ws.Range("I45").Value = "18.0 mm"
$ws.Range("J45").Value = "1017.9 hPa"
$ws.Range("E46").Value = "2026-02-16 21:20:24"
$ws.Range("J46").Value = "1016.8 hPa"

# Percent-formatted text cells need special handling: assigning a
# "NN%" string straight to .Value gets auto-parsed into a numeric
# percentage by the input-inference logic, which would also rewrite
# the cell number format. These columns store plain text in this
# sheet, so build the literal text via a formula, then collapse the
# formula down to a static value (Copy + PasteSpecial values-only) so
# the stored cell keeps its original text type and style.
$ws.Range("H4").Formula = "=""61%"""
$ws.Range("H4").Copy() | Out-Null
$ws.Range("H4").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("H10").Formula = "=""77%"""
$ws.Range("H10").Copy() | Out-Null
$ws.Range("H10").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("H15").Formula = "=""66%"""
$ws.Range("H15").Copy() | Out-Null
$ws.Range("H15").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("H18").Formula = "=""75%"""
$ws.Range("H18").Copy() | Out-Null
$ws.Range("H18").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("H30").Formula = "=""69%"""
$ws.Range("H30").Copy() | Out-Null
$ws.Range("H30").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("H35").Formula = "=""74%"""
$ws.Range("H35").Copy() | Out-Null
$ws.Range("H35").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("H37").Formula = "=""83%"""
$ws.Range("H37").Copy() | Out-Null
$ws.Range("H37").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("H38").Formula = "=""71%"""
$ws.Range("H38").Copy() | Out-Null
$ws.Range("H38").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("H45").Formula = "=""96%"""
$ws.Range("H45").Copy() | Out-Null
$ws.Range("H45").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = 0
